# Generate Report for Handback
# Updates the Overview / zh-cn / de-de sheets with the handback report:
#  - "In Translation" -> "Handed back: in sync with en-US" (status/summary cells)
#  - "Latest Target File" gets the source markdown file name (as a hyperlink)
#  - "Latest Handback File" gets the generated locale xlf filename
#  - "Latest Handback DateTime" gets the handback timestamp

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $statusText
$ovw.Range("F2").Value = $statusText
$ovw.Range("E3").Value = $statusText
$ovw.Range("F3").Value = $statusText

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Range("I2").Value = "944f5352-7142-419b-8ff0-2e96da6128dd.md"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/749336c55252609c5885a28d1e581e7701bb6bf2/e2e/944f5352-7142-419b-8ff0-2e96da6128dd.md", "", "", "944f5352-7142-419b-8ff0-2e96da6128dd.md") | Out-Null
$zh.Range("J2").Value = "944f5352-7142-419b-8ff0-2e96da6128dd.11b92db4792c49cccf911972f250efeb4353df4b.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-25 14:23:19"

$zh.Range("I3").Value = "acf06688-41a3-4799-9226-34af7ba81abe.md"
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/749336c55252609c5885a28d1e581e7701bb6bf2/e2e/acf06688-41a3-4799-9226-34af7ba81abe.md", "", "", "acf06688-41a3-4799-9226-34af7ba81abe.md") | Out-Null
$zh.Range("J3").Value = "acf06688-41a3-4799-9226-34af7ba81abe.690551d0200f2c8a192c6b86019f1e70f235ded6.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-25 14:23:19"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Range("I2").Value = "944f5352-7142-419b-8ff0-2e96da6128dd.md"
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/749336c55252609c5885a28d1e581e7701bb6bf2/e2e/944f5352-7142-419b-8ff0-2e96da6128dd.md", "", "", "944f5352-7142-419b-8ff0-2e96da6128dd.md") | Out-Null
$de.Range("J2").Value = "944f5352-7142-419b-8ff0-2e96da6128dd.11b92db4792c49cccf911972f250efeb4353df4b.de-de.xlf"
$de.Range("K2").Value = "2016-08-25 14:23:26"

$de.Range("I3").Value = "acf06688-41a3-4799-9226-34af7ba81abe.md"
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/749336c55252609c5885a28d1e581e7701bb6bf2/e2e/acf06688-41a3-4799-9226-34af7ba81abe.md", "", "", "acf06688-41a3-4799-9226-34af7ba81abe.md") | Out-Null
$de.Range("J3").Value = "acf06688-41a3-4799-9226-34af7ba81abe.690551d0200f2c8a192c6b86019f1e70f235ded6.de-de.xlf"
$de.Range("K3").Value = "2016-08-25 14:23:26"
